# DT added label option
# The backward-elimination report was re-run; every sheet's OLS summary
# (held in cell B2) gets a refreshed "Date:" / "Time:" stamp while every
# other figure in the report stays identical.

$wb = $excel.ActiveWorkbook

$oldDate = "Date:                Sun, 05 Jan 2020"
$newDate = "Date:                Wed, 08 Jan 2020"
$oldTime = "Time:                        21:22:49"
$newTime = "Time:                        19:07:52"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $summary = $cell.Value2
    if ($summary -ne $null -and $summary.Contains("Date:")) {
        $summary = $summary.Replace($oldDate, $newDate)
        $summary = $summary.Replace($oldTime, $newTime)
        $cell.Value2 = $summary
    }
}
